$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: shuffle "Single Var Calc" / "EdX" across the week
$ws.Range("C12").Value = "Single Var Calc"
$ws.Range("D12").Value = "Single Var Calc"
$ws.Range("E12").Value = "EdX"
$ws.Range("F12").Value = "EdX"
$ws.Range("G12").Value = "EdX"
$ws.Range("H12").Value = "Single Var Calc"
$ws.Range("I12").Value = "EdX"

# Row 15: swap of "EdX" / "Single Var Calc" (previously "Read SCIP")
$ws.Range("C15").Value = "EdX"
$ws.Range("D15").Value = "EdX"
$ws.Range("E15").Value = "Single Var Calc"
$ws.Range("F15").Value = "Single Var Calc"
$ws.Range("G15").Value = "Single Var Calc"
$ws.Range("H15").Value = "EdX"
$ws.Range("I15").Value = "Single Var Calc"

# Row 16: "Personal Projects" row cleared out (moved down to row 20)
$ws.Range("C16:I16").ClearContents()

# Row 19: "Lunch" stays the same text, re-set for safety
$ws.Range("C19:I19").Value = "Lunch"

# Row 20: "Personal Projects" moved here from row 16
$ws.Range("C20").Value = "Personal Projects"
$ws.Range("D20").Value = "Personal Projects"
$ws.Range("E20").Value = "Personal Projects"
$ws.Range("F20").Value = "Personal projects"
$ws.Range("G20").Value = "Personal Projects"
$ws.Range("H20").Value = "Personal Projects"
$ws.Range("I20").Value = "Personal Projects"

# Row 27: newly filled in with "Read SCIP"
$ws.Range("C27:I27").Value = "Read SCIP"
# F27 ends up with its own distinct (but visually identical) cell format
$ws.Range("F27").Font.ThemeColor = 5

# Row 35: "TEST" removed
$ws.Range("C35").ClearContents()

# Row 40: newly filled in with "BED" (set before row 38 so shared-string order matches)
$ws.Range("C40:I40").Value = "BED"

# Row 38: newly filled in with "Leisure Reading"
$ws.Range("C38:I38").Value = "Leisure Reading"

# Update the view state to match (scrolled/selected area)
$ws.Range("H41").Select()
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
